$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) cells to Text format so numeric-looking strings
# like "46.078.57" or "1.00" are preserved exactly as text, not coerced to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.078.57"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.584.96"
$ws.Range("E3").Value = "  +8.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.50"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.89"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.593"
$ws.Range("E7").Value = "  +4.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.575"
$ws.Range("E9").Value = "  +13.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.45"
$ws.Range("E10").Value = "  +11.64%  "
$ws.Range("E11").Value = "  +6.21%  "
$ws.Range("E12").Value = "  +14.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.978.42"
$ws.Range("E13").Value = "  +8.58%  "
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.605.78"
$ws.Range("E15").Value = "  +8.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.897"
$ws.Range("E16").Value = "  +8.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.77"
$ws.Range("E17").Value = "  +7.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.206.78"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000101"
$ws.Range("E19").Value = "  +6.28%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.02"
$ws.Range("E20").Value = "  +3.01%  "
$ws.Range("E21").Value = "  +9.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.83"
$ws.Range("E22").Value = "  +6.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.27"
$ws.Range("E23").Value = "  +4.34%  "
$ws.Range("E24").Value = "  +6.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").Value = "  +14.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.89"
$ws.Range("E26").Value = "  +33.40%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.40"
$ws.Range("E28").Value = "  +7.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.45"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("E31").Value = "  +9.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.69"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.27"
$ws.Range("E34").Value = "  +18.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "152.41"
$ws.Range("E35").Value = "  +4.27%  "
$ws.Range("E36").Value = "  +6.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.116"
$ws.Range("E37").Value = "  +3.45%  "
$ws.Range("E38").Value = "  +4.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.12"
$ws.Range("E39").Value = "  +8.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.17"
$ws.Range("E40").Value = "  +7.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.57"
$ws.Range("E41").Value = "  +11.89%  "
$ws.Range("E42").Value = "  +7.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.061.47"
$ws.Range("E43").Value = "  +6.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.73"
$ws.Range("E44").Value = "  +38.21%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.25"
$ws.Range("E47").Value = "  +9.11%  "
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "108.30"
$ws.Range("E49").Value = "  +9.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.199"
$ws.Range("E50").Value = "  +7.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.837.65"
$ws.Range("E51").Value = "  +8.60%  "
